$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header column F1 = "EDAM_DEF", matching the existing header
# formatting (bold, bordered, centered) already used by B1:E1.
$ws.Range("F1").Value = "EDAM_DEF"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
